$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.207.52"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "2.353.50"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("E5").Value = "  +3.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.98"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.31"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +6.71%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.566"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +23.27%  "
$ws.Range("E10").Value = "  +4.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.72"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +21.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.36"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +18.49%  "
$ws.Range("E13").Value = "  +1.95%  "
$ws.Range("D14").Value = "2.704.24"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.79"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +7.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.907"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.99%  "
$ws.Range("D17").Value = "2.362.83"
$ws.Range("E17").Value = "  -2.84%  "
$ws.Range("D18").Value = "44.390.48"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000103"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.91%  "
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "77.90"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +5.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("E24").Value = "  -5.15%  "
$ws.Range("E25").Value = "  +3.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.66"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.24%  "
$ws.Range("E27").Value = "  +1.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.52"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.71"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("E30").Value = "  +3.29%  "
$ws.Range("E31").Value = "  +3.78%  "
$ws.Range("E32").Value = "  +4.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.37"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +8.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0761"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +10.56%  "
$ws.Range("E35").Value = "  +4.87%  "
$ws.Range("E36").Value = "  +7.12%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.54"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("E39").Value = "  +7.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.34"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.96"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1000"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.190"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +14.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.62"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.49%  "
$ws.Range("E47").Value = "  +10.00%  "
$ws.Range("E48").Value = "  -2.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.48"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.41%  "
$ws.Range("D50").Value = "1.450.64"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000207"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.92%  "
